$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.442.65"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.945.08"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.23"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.31"
$ws.Range("E6").Value = "  +5.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.943.24"
$ws.Range("E8").Value = "  +2.24%  "
$ws.Range("E9").Value = "  +3.18%  "
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("E11").Value = "  +9.61%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +7.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.23"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.434.19"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.482.50"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.945.19"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.76"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.05"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.91"
$ws.Range("E26").Value = "  +5.26%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E29").Value = "  +7.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.57"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("E31").Value = "  +4.35%  "
$ws.Range("E32").Value = "  +18.43%  "
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.15"
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.989"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("E38").Value = "  +7.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.72"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.34"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  +4.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.63"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "135.04"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.685.34"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "354.99"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.46"
$ws.Range("E51").Value = "  +0.34%  "
